$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 13 (shifts old rows 13-23 down to 14-24).
#    This reproduces the "Docentes responsaveis -> (new row with professor name)" insertion
$ws.Rows.Item(13).Insert()

# 2) Clean up the newly-inserted row 13: no A13 cell; B13/C13 hold the professor-name value
#    that used to be mis-placed under "Objetivos:" (row 10).
$ws.Range("D1").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").ClearContents()

$ws.Range("B13").Value = '5840535 - Messias Borges Silva'
$ws.Range("B2").Copy()
$ws.Range("B13").PasteSpecial(-4122)

$ws.Range("C13").Value = '5840535 - Messias Borges Silva'
$ws.Range("C2").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# 3) Fix up the content of the rows whose values were wrong / placeholders in the source file.

# Row 10 ("Objetivos:") previously held the misplaced professor name -> real PT objectives text
$ws.Range("B10").Value = 'Introduzir os alunos nos conceitos técnicos fundamentais de um curso de Engenharia de Produção, tendo em vista a sua formação generalista voltada para os mais diversos tipos de sistemas de produção.'
$ws.Range("C10").Value = 'Introduzir os alunos nos conceitos técnicos fundamentais de um curso de Engenharia de Produção, tendo em vista a sua formação generalista voltada para os mais diversos tipos de sistemas de produção.'

# Row 14 ("Programa resumido:") previously held "Semestral" -> real PT short-syllabus text
$ws.Range("B14").Value = '1 – Introdução aos Sistemas Produtivos2 – Papel Estratégico da Produção3 – Estratégia de Produção4 – Projeto em Gestão de Produção5 – Projeto de Produtos e Serviços6 – Projeto da Rede de Operações Produtivas7 – Arranjo Físico e Fluxo'
$ws.Range("C14").Value = '1 – Introdução aos Sistemas Produtivos2 – Papel Estratégico da Produção3 – Estratégia de Produção4 – Projeto em Gestão de Produção5 – Projeto de Produtos e Serviços6 – Projeto da Rede de Operações Produtivas7 – Arranjo Físico e Fluxo'

# Row 16 ("Programa:") previously held "01/01/2018" -> real PT full-syllabus text
$ws.Range("B16").Value = '1 – Introdução aos Sistemas ProdutivosProdução na Organização. Inputs, Processos de Transformação e Outputs. Tipos de Operações de Produção. Atividades da administração da produção.2 – Papel Estratégico da ProduçãoPapel da função produção. Objetivos de desempenho. 3 – Estratégia de ProduçãoIntrodução. Prioridade de objetivos de desempenho. Áreas de decisão da estratégia de operações.4 – Projeto em Gestão de ProduçãoDefinição de projeto. Principais aspectos de um projeto. Tipos de processos em manufatura e serviços. 5 – Projeto de Produtos e ServiçosIntrodução. Geração de conceito. Triagem de conceito. Avaliação e melhoria do projeto. Protótipo e projeto final.6 – Projeto da Rede de Operações ProdutivasPerspectiva da rede. Integração Vertical. Localização da capacidade. Gestão da capacidade produtiva de longo prazo.7 – Arranjo Físico e FluxoProcedimento de Arranjo Físico. Tipos básicos de arranjo físico. Projeto de arranjo físico.'
$ws.Range("C16").Value = '1 – Introdução aos Sistemas ProdutivosProdução na Organização. Inputs, Processos de Transformação e Outputs. Tipos de Operações de Produção. Atividades da administração da produção.2 – Papel Estratégico da ProduçãoPapel da função produção. Objetivos de desempenho. 3 – Estratégia de ProduçãoIntrodução. Prioridade de objetivos de desempenho. Áreas de decisão da estratégia de operações.4 – Projeto em Gestão de ProduçãoDefinição de projeto. Principais aspectos de um projeto. Tipos de processos em manufatura e serviços. 5 – Projeto de Produtos e ServiçosIntrodução. Geração de conceito. Triagem de conceito. Avaliação e melhoria do projeto. Protótipo e projeto final.6 – Projeto da Rede de Operações ProdutivasPerspectiva da rede. Integração Vertical. Localização da capacidade. Gestão da capacidade produtiva de longo prazo.7 – Arranjo Físico e FluxoProcedimento de Arranjo Físico. Tipos básicos de arranjo físico. Projeto de arranjo físico.'

# Rows 19-22 ("Metodo:", "Criterio:", "Norma de recuperacao:", "Bibliografia:") were each
# holding the VALUE belonging to the next row (off-by-one) -> shift them back into place
# and give "Bibliografia:" its real (new) value.
$ws.Range("B19").Value = 'Aulas Expositivas; trabalhos e seminários.'
$ws.Range("C19").Value = 'Aulas Expositivas; trabalhos e seminários.'

$ws.Range("B20").Value = 'MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.'
$ws.Range("C20").Value = 'MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.'

$ws.Range("B21").Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação.'
$ws.Range("C21").Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação.'

$ws.Range("B22").Value = 'SLACK, N. et al. Administração da produção. São Paulo: Atlas, 2002. 
Textos complementares serão usados durante o curso.'
$ws.Range("C22").Value = 'SLACK, N. et al. Administração da produção. São Paulo: Atlas, 2002. 
Textos complementares serão usados durante o curso.'

# 4) Split the merged column A/B width definition: column A keeps its own width entry and
#    column B keeps its already-distinct width/style entry.
$ws.Columns.Item(2).ColumnWidth = 60.7109375

